$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top
$ws.Rows.Item(1).Insert()

# Set the value of the new A1 cell
$ws.Range("A1").Value = "QFR"

# Update selection to A2
$ws.Range("A2").Select()
